{"js": "// Ultima IV Dokumentation / Referenz.docx \u2014 proofing pass.\n//\n// The author accepted a round of Word grammar/hyphenation suggestions.\n// Almost all of the underlying XML churn (stale <w:proofErr/> markers\n// being cleared, runs being re-split at new soft line/page-break points,\n// a new `oel` namespace + list `durableId`s stamped by a newer Word build)\n// is not reachable through the Word automation surface \u2014 Word itself\n// regenerates those on its internal proofing/pagination passes and does\n// not expose them as editable objects. What *is* reachable, and what\n// actually changes the visible document, are:\n//   1. A stray double space fixed (\"auf  deinen\" -> \"auf deinen\").\n//   2. A grammar fix + clause reorder in the \"GEBEN\" paragraph.\n//   3. A word choice swap (\"einen Umhang\" -> \"eine Robe\").\n//   4. A few soft hyphens added/removed at new hyphenation points.\n//\n// Soft hyphens aren't literal text, but this host's Range.text surfaces\n// them as U+001F (matching the VBA/COM Chr(31) convention), and feeding\n// that same character back into insertText() round-trips it to a proper\n// <w:softHyphen/> run on save \u2014 so we use it below for both search\n// (via wildcards, since a literal search string can't embed it) and\n// insertion.\n\nconst SOFT_HYPHEN = \"\\u001F\";\nconst body = context.document.body;\n\nasync function replaceOnce(searchText, replacement, options) {\n  const results = body.search(searchText, Object.assign({ matchCase: true }, options));\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match for \" + JSON.stringify(searchText) +\n      \" but found \" + results.items.length\n    );\n  }\n\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1. \"Zu Fu\u00df [G]ehen\" paragraph: collapse the double space before\n//    \"deinen\" and add the new hyphenation point inside it.\nawait replaceOnce(\"auf  deinen\", \"auf dei\" + SOFT_HYPHEN + \"nen\");\n\n// 2. \"[E]nde & Speichern\" paragraph: new hyphenation point in \"Reise\".\nawait replaceOnce(\"deine Reise beenden\", \"deine Rei\" + SOFT_HYPHEN + \"se beenden\");\n\n// 3. \"W\u00e4hrend deiner Gespr\u00e4che ...\" paragraph: de-capitalize \"Du\" and\n//    move the \"Du kannst das tun\" clause from the front to the end.\nawait replaceOnce(\"k\u00f6nntest Du den\", \"k\u00f6nntest du den\");\nawait replaceOnce(\"zeigen. Du kannst das tun, indem du\", \"zeigen. Indem du\");\nawait replaceOnce(\"sagst.\", \"sagst, kannst du das tun.\");\n\n// 4. \"Einige der Leute ...\" / Iolo paragraph: word choice change.\nawait replaceOnce(\"einen Umhang\", \"eine Robe\");\n\n// 5. \"Viele der Questen von Ultima IV ...\" paragraph: remove the stale\n//    hyphenation point inside \"Questen\" (wildcard search since the\n//    existing soft hyphen sits inside the literal match).\nawait replaceOnce(\"Viele der Q*sten von Ultima\", \"Viele der Questen von Ultima\", { matchWildcards: true });\n\n// 6. Piratenschiff paragraph: new hyphenation point in \"Rumpfst\u00e4rke\".\nawait replaceOnce(\n  \"dann w\u00fcrde die Rumpfst\u00e4rke des Schiffes\",\n  \"dann w\u00fcrde die Rumpfst\u00e4r\" + SOFT_HYPHEN + \"ke des Schiffes\"\n);\n", "ps1": "# Ultima IV Dokumentation / Referenz.docx - proofing pass.\n#\n# The author accepted a round of Word grammar/hyphenation suggestions.\n# Almost all of the underlying XML churn (stale proofErr markers being\n# cleared, runs being re-split at new soft line/page-break points, a new\n# `oel` namespace + list durableIds stamped by a newer Word build) is not\n# reachable through Word automation - Word regenerates those internally\n# during its own proofing/pagination passes and doesn't expose them as\n# editable objects. What *is* reachable, and what actually changes the\n# visible document, are:\n#   1. A stray double space fixed (\"auf  deinen\" -> \"auf deinen\").\n#   2. A grammar fix + clause reorder in the \"GEBEN\" paragraph.\n#   3. A word choice swap (\"einen Umhang\" -> \"eine Robe\").\n#   4. A few soft hyphens added/removed at new hyphenation points.\n#\n# Soft hyphens aren't literal text, but this host surfaces them through\n# Range.Text as Chr(31) (the usual VBA/COM convention), and writing that\n# same character back round-trips it to a proper <w:softHyphen/> run on\n# save - so we use it below both to match existing soft hyphens (via\n# wildcard search, since a literal search string can't embed one) and to\n# insert new ones.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Once {\n    param(\n        [string]$SearchText,\n        [string]$Replacement,\n        [bool]$Wildcards = $false\n    )\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $SearchText\n    $find.MatchWildcards = $Wildcards\n    $find.MatchCase = $true\n    $find.Forward = $true\n\n    $found = $find.Execute()\n    if (-not $found) {\n        throw \"Replace-Once: no match for '$SearchText'\"\n    }\n\n    # $rng now covers exactly the matched text; replacing its .Text\n    # substitutes the run(s) in place and leaves surrounding formatting\n    # (e.g. the bold \"GEBEN\" run) untouched.\n    $rng.Text = $Replacement\n}\n\n$SoftHyphen = [char]31\n\n# 1. \"Zu Fuss [G]ehen\" paragraph: collapse the double space before\n#    \"deinen\" and add the new hyphenation point inside it.\nReplace-Once \"auf  deinen\" (\"auf dei\" + $SoftHyphen + \"nen\")\n\n# 2. \"[E]nde & Speichern\" paragraph: new hyphenation point in \"Reise\".\nReplace-Once \"deine Reise beenden\" (\"deine Rei\" + $SoftHyphen + \"se beenden\")\n\n# 3. \"Waehrend deiner Gespraeche ...\" paragraph: de-capitalize \"Du\" and\n#    move the \"Du kannst das tun\" clause from the front to the end.\nReplace-Once \"k\u00f6nntest Du den\" \"k\u00f6nntest du den\"\nReplace-Once \"zeigen. Du kannst das tun, indem du\" \"zeigen. Indem du\"\nReplace-Once \"sagst.\" \"sagst, kannst du das tun.\"\n\n# 4. Iolo paragraph: word choice change.\nReplace-Once \"einen Umhang\" \"eine Robe\"\n\n# 5. \"Viele der Questen von Ultima IV ...\" paragraph: remove the stale\n#    hyphenation point inside \"Questen\" (wildcard search since the\n#    existing soft hyphen sits inside the literal match).\nReplace-Once \"Viele der Q*sten von Ultima\" \"Viele der Questen von Ultima\" $true\n\n# 6. Piratenschiff paragraph: new hyphenation point in \"Rumpfstaerke\".\nReplace-Once \"dann w\u00fcrde die Rumpfst\u00e4rke des Schiffes\" (\"dann w\u00fcrde die Rumpfst\u00e4r\" + $SoftHyphen + \"ke des Schiffes\")\n"}
